$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Optimal_length" column (C),
# pushing it to D. The new column C will hold the "corrected" (upravene)
# optimal length values, while D keeps the original ("wu") values.
$ws.Range("C1").EntireColumn.Insert()

# Headers
$ws.Range("C1").Formula = '="Optimal_length_upravene"'
$ws.Range("C1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4163) | Out-Null

$ws.Range("D1").Formula = '="Optimal_length_wu"'
$ws.Range("D1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4163) | Out-Null

# Column C mirrors column D for every row except SRA1 (row 3), which gets
# a newly fitted value. Values must stay text (shared-string) typed, same
# as the rest of the sheet, so they are produced as formulas and then
# converted to static values via copy/paste-values (keeps default style).
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 4).Copy() | Out-Null
    $ws.Cells.Item($r, 3).PasteSpecial(-4163) | Out-Null
}

$ws.Range("C3").Formula = '="0.1252089719457094"'
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = 0

# Column widths (closest values the runtime's character-width quantization
# can reach to the saved widths of 25.7109375 / 19.28515625)
$ws.Range("C1").EntireColumn.ColumnWidth = 24.833333333333336
$ws.Range("D1").EntireColumn.ColumnWidth = 18.5

# Restore the selection as saved in the workbook
$ws.Range("C4").Select()
